# library_3357.xlsx — remove the duplicate fastq-library entry.
# The original row 2 (s2cDNADate 11.14.18 / index1Sequence TGAGGTTATC,
# sample #1) was a duplicate of the row that is now #7 (same date/index),
# so the whole row is deleted and everything below shifts up one row,
# matching the commit "deleted duplicate fastq files in 3357 and 3275 ...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire second row (the duplicate entry); remaining rows shift up.
$ws.Rows(2).Delete()

# Leave the same kind of "entire row" selection Excel leaves behind after
# a row deletion (row 2, now holding what used to be row 3's data).
$ws.Rows(2).Select() | Out-Null
